$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.010.92'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').Value = '2.231.09'
$ws.Range('E3').Value = '  -6.50%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.56'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -5.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '80.06'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -9.16%  '
$ws.Range('E7').Value = '  -4.42%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -7.15%  '
$ws.Range('E10').Value = '  -6.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '27.91'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -10.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.38'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -12.60%  '
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').Value = '2.573.32'
$ws.Range('E14').Value = '  -6.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.10'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -7.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.04'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -7.73%  '
$ws.Range('D17').Value = '2.234.67'
$ws.Range('E17').Value = '  -6.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.714'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -6.32%  '
$ws.Range('D19').Value = '38.934.06'
$ws.Range('E19').Value = '  -4.43%  '
$ws.Range('D20').Value = '0.0₃0857'
$ws.Range('E20').Value = '  -6.16%  '
$ws.Range('E21').Value = '  -7.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.83'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.80'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -9.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '224.89'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.44%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.39'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -10.27%  '
$ws.Range('E27').Value = '  -6.47%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.14'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.56%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.13'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.88'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '148.94'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.02'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -8.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('E34').Value = '  -9.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.33'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('E36').Value = '  -6.81%  '
$ws.Range('E37').Value = '  -4.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.65'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.21%  '
$ws.Range('E39').Value = '  -5.26%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -10.19%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.60'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -8.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.63'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.08%  '
$ws.Range('D43').Value = '1.903.89'
$ws.Range('E43').Value = '  -3.00%  '
$ws.Range('E44').Value = '  -9.18%  '
$ws.Range('E45').Value = '  -6.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.95'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.57%  '
$ws.Range('E48').Value = '  -10.76%  '
$ws.Range('D49').Value = '2.439.20'
$ws.Range('E49').Value = '  -6.55%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.32'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.43%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '87.37'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -7.03%  '
